$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.925.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.75%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.909.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.73%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'319.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.68%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E7").Value = "'  -2.40%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4052"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.97%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'41.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.89%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'24.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.88%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.907.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.12%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.392"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.43%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.67%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'92.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.81%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001098"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.97%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06505"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.49%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'18.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.74%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.13%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.935"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.12%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'29.959.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.68%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.11%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.201"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.57%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'22.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.28%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.124.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.18%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'161.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.299"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.56%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'128.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.06%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.59%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.85%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.929"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.98%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.804"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.35%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.401"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.90%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.02439"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.96%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.06408"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.17%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.82%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.194"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.80%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'8.707"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.20%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6460"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'11.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.42%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.208"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.19%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +8.41%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'13.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.88%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.6031"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.17%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.93%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'122.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.67%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.32%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'78.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.33%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.126"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.55%  "
$ws.Range("E51").Style = "Normal"
